$wb = $excel.ActiveWorkbook

# --- Sheet "Variables": lower-case three device codes (rows 10-12, column B) ---
$wsVariables = $wb.Worksheets.Item("Variables")
$wsVariables.Range('B10').Value = 'zt3'
$wsVariables.Range('B11').Value = 'zr3'
$wsVariables.Range('B12').Value = 'pf3'

# --- Sheet "Categories": rebuild rows 40-71 with the new category table ---
$wsCategories = $wb.Worksheets.Item("Categories")

$wsCategories.Cells.Item(40, 1).Value = 'casemi_fup5'
$wsCategories.Cells.Item(40, 2).Value = 1
$wsCategories.Cells.Item(40, 3).Value = 'prevalent'

$wsCategories.Cells.Item(41, 1).Value = 'casemi_fup5'
$wsCategories.Cells.Item(41, 2).Value = 2
$wsCategories.Cells.Item(41, 3).Value = 'incident (verif.)'

$wsCategories.Cells.Item(42, 1).Value = 'casemi_fup5'
$wsCategories.Cells.Item(42, 2).Value = 5
$wsCategories.Cells.Item(42, 3).Value = 'I252 (old MI, unknown date of occurence)'

$wsCategories.Cells.Item(43, 1).Value = 'casemi_fup5'
$wsCategories.Cells.Item(43, 2).Value = 9
$wsCategories.Cells.Item(43, 3).Value = 'incident (not verif.)'

$wsCategories.Cells.Item(44, 1).Value = 'casestroke_fup5'
$wsCategories.Cells.Item(44, 2).Value = 0
$wsCategories.Cells.Item(44, 3).Value = 'not diseased'

$wsCategories.Cells.Item(45, 1).Value = 'casestroke_fup5'
$wsCategories.Cells.Item(45, 2).Value = 1
$wsCategories.Cells.Item(45, 3).Value = 'prevalent'

$wsCategories.Cells.Item(46, 1).Value = 'casestroke_fup5'
$wsCategories.Cells.Item(46, 2).Value = 2
$wsCategories.Cells.Item(46, 3).Value = 'incident (verif.)'

$wsCategories.Cells.Item(47, 1).Value = 'casestroke_fup5'
$wsCategories.Cells.Item(47, 2).Value = 6
$wsCategories.Cells.Item(47, 3).Value = 'I64Y old stroke, date of diagnosis unknown'

$wsCategories.Cells.Item(48, 1).Value = 'casestroke_fup5'
$wsCategories.Cells.Item(48, 2).Value = 9
$wsCategories.Cells.Item(48, 3).Value = 'incident (not verif.)'

$wsCategories.Cells.Item(49, 1).Value = 'caseI63_fup5'
$wsCategories.Cells.Item(49, 2).Value = 0
$wsCategories.Cells.Item(49, 3).Value = 'No'

$wsCategories.Cells.Item(50, 1).Value = 'caseI63_fup5'
$wsCategories.Cells.Item(50, 2).Value = 1
$wsCategories.Cells.Item(50, 3).Value = 'Yes'

$wsCategories.Cells.Item(51, 1).Value = 'caseI61_fup5'
$wsCategories.Cells.Item(51, 2).Value = 0
$wsCategories.Cells.Item(51, 3).Value = 'No'

$wsCategories.Cells.Item(52, 1).Value = 'caseI61_fup5'
$wsCategories.Cells.Item(52, 2).Value = 1
$wsCategories.Cells.Item(52, 3).Value = 'Yes'

$wsCategories.Cells.Item(53, 1).Value = 'casehyp_fup5'
$wsCategories.Cells.Item(53, 2).Value = 0
$wsCategories.Cells.Item(53, 3).Value = 'not diseased'

$wsCategories.Cells.Item(54, 1).Value = 'casehyp_fup5'
$wsCategories.Cells.Item(54, 2).Value = 1
$wsCategories.Cells.Item(54, 3).Value = 'prevalent'

$wsCategories.Cells.Item(55, 1).Value = 'casehyp_fup5'
$wsCategories.Cells.Item(55, 2).Value = 2
$wsCategories.Cells.Item(55, 3).Value = 'incident (verif.)'

$wsCategories.Cells.Item(56, 1).Value = 'casehyp_fup5'
$wsCategories.Cells.Item(56, 2).Value = 3
$wsCategories.Cells.Item(56, 3).Value = 'incident I15'

$wsCategories.Cells.Item(57, 1).Value = 'casehyp_fup5'
$wsCategories.Cells.Item(57, 2).Value = 9
$wsCategories.Cells.Item(57, 3).Value = 'incident (not verif.)'

$wsCategories.Cells.Item(58, 1).Value = 'casehf_fup5'
$wsCategories.Cells.Item(58, 2).Value = 0
$wsCategories.Cells.Item(58, 3).Value = 'not diseased'

$wsCategories.Cells.Item(59, 1).Value = 'casehf_fup5'
$wsCategories.Cells.Item(59, 2).Value = 1
$wsCategories.Cells.Item(59, 3).Value = 'prevalent'

$wsCategories.Cells.Item(60, 1).Value = 'casehf_fup5'
$wsCategories.Cells.Item(60, 2).Value = 2
$wsCategories.Cells.Item(60, 3).Value = 'incident (verif.)'

$wsCategories.Cells.Item(61, 1).Value = 'casehf_fup5'
$wsCategories.Cells.Item(61, 2).Value = 9
$wsCategories.Cells.Item(61, 3).Value = 'incident (not verif.)'

$wsCategories.Cells.Item(62, 1).Value = 'casediab_fup5'
$wsCategories.Cells.Item(62, 2).Value = 0
$wsCategories.Cells.Item(62, 3).Value = 'not diseased'

$wsCategories.Cells.Item(63, 1).Value = 'casediab_fup5'
$wsCategories.Cells.Item(63, 2).Value = 1
$wsCategories.Cells.Item(63, 3).Value = 'prevalent'

$wsCategories.Cells.Item(64, 1).Value = 'casediab_fup5'
$wsCategories.Cells.Item(64, 2).Value = 2
$wsCategories.Cells.Item(64, 3).Value = 'incident (verif.)'

$wsCategories.Cells.Item(65, 1).Value = 'casediab_fup5'
$wsCategories.Cells.Item(65, 2).Value = 4
$wsCategories.Cells.Item(65, 3).Value = 'inc. Diabetes (other types)'

$wsCategories.Cells.Item(66, 1).Value = 'casediab_fup5'
$wsCategories.Cells.Item(66, 2).Value = 9
$wsCategories.Cells.Item(66, 3).Value = 'incident (not verif.)'

$wsCategories.Cells.Item(67, 1).Value = 'inccanc_fup5'
$wsCategories.Cells.Item(67, 2).Value = 0
$wsCategories.Cells.Item(67, 3).Value = 'No'

$wsCategories.Cells.Item(68, 1).Value = 'inccanc_fup5'
$wsCategories.Cells.Item(68, 2).Value = 1
$wsCategories.Cells.Item(68, 3).Value = 'Yes'

$wsCategories.Cells.Item(69, 1).Value = 'vitstat5'
$wsCategories.Cells.Item(69, 2).Value = 1
$wsCategories.Cells.Item(69, 3).Value = 'alive'

$wsCategories.Cells.Item(70, 1).Value = 'vitstat5'
$wsCategories.Cells.Item(70, 2).Value = 2
$wsCategories.Cells.Item(70, 3).Value = 'dead'

$wsCategories.Cells.Item(71, 1).Value = 'vitstat5'
$wsCategories.Cells.Item(71, 2).Value = 6
$wsCategories.Cells.Item(71, 3).Value = 'dropped out'

